$wb = $excel.ActiveWorkbook

# --- Sheet1 "TestSteps": row 5's KeyInData changes to "getData=Profile" ---
$ws1 = $wb.Worksheets.Item("TestSteps")
$ws1.Range("C5").Value = "getData=Profile"

# --- Sheet2 "TestData": add new "Profile" column E ---
$ws2 = $wb.Worksheets.Item("TestData")

# E1 header - copy style from D1 (header style), then set its value to "Profile"
$ws2.Range("D1").Copy($ws2.Range("E1"))
$ws2.Range("E1").Value = "Profile"

# E2 - copy style from D2 (data style), then set value "n/a"
$ws2.Range("D2").Copy($ws2.Range("E2"))
$ws2.Range("E2").Value = "n/a"

# D3 gets "n/a" (already has correct style, just set value)
$ws2.Range("D3").Value = "n/a"

# E3 - copy style from D2 (empty data style cell), stays empty
$ws2.Range("D2").Copy($ws2.Range("E3"))
$ws2.Range("E3").ClearContents()

# Update the remembered selections
$ws2.Range("E8").Select()
$ws1.Range("C14").Select()
